$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.395.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.40%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.723.17'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.39%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.84%  '

# Row 6
$ws.Range("E6").Value = '  -0.02%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4865'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.55%  '

# Row 8
$ws.Range("E8").Value = '  -3.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06192'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.49%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.715.16'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.89%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06974'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.25%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.47'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.05%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.533'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.56%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5966'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.44%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.18'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.25%  '

# Row 16
$ws.Range("E16").Value = '  -0.04%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.401.88'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.34%  '

# Row 18
$ws.Range("E18").Value = '  +0.01%  '

# Row 19
$ws.Range("E19").Value = '  +0.04%  '

# Row 20
$ws.Range("E20").Value = '  -2.11%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.951.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.05%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.442'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.53%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.471'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.59%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.099'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.02%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '137.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.16%  '

# Row 26
$ws.Range("E26").Value = '  -1.16%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.400'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.20%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '106.54'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.63%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.723'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.07%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.923'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.53%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08001'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.13%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.666'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.75%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04496'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.49%  '

# Row 34
$ws.Range("E34").Value = '  -0.49%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9978'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.75%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6237'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.69%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9334'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.42%  '

# Row 38
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.387'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.05%  '

# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.947'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.37%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9997'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.41%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01474'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.98%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '100.02'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.59%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.332'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.71%  '

# Row 44
$ws.Range("E44").Value = '  -1.53%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.869'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.19%  '

# Row 46
$ws.Range("E46").Value = '  -1.70%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05365'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.32%  '

# Row 48
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.713'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.58%  '

# Row 49
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.09'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.62%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.227'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.91%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '50.87'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.29%  '
